$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date-number-format style (s="1", numFmtId 14) from A38 down through A39:A55
# via a format-only paste, so no new style/numFmt entries are created.
$ws.Range("A38").Copy()
$ws.Range("A39:A55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row (A=date serial, B=URL) tuples, listed in the exact chronological order the
# original author entered them (rows 42/43 were filled out of sequence), so the
# resulting shared-string table indices line up with the target workbook.
$entries = @(
  @{Row=39; Date=44291; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/1ab5a678c34cce2fa4e095f158cc1b152bac64d0/states.json"},
  @{Row=40; Date=44292; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/7daa747efd2b02f3e1bd61c2d9844fe7761bb88d/states.json"},
  @{Row=41; Date=44293; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/59dc3c257afd9d129702f0fa3cac73945923eb75/states.json"},
  @{Row=43; Date=44295; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8f182bd89dfc7efac6033cdff502dead207a5c9a/states.json"},
  @{Row=42; Date=44294; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/2b88b4a7afeb95cc81c4e4bf834a3a8981c94285/states.json"},
  @{Row=44; Date=44296; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/b46a9e4137a23a86b08770d3c737ee9ae84a051a/states.json"},
  @{Row=45; Date=44297; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/4d055f8550e1c5635d7425632b7be582edee6377/states.json"},
  @{Row=46; Date=44298; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/06ed24aa6756de13a6d1d0c283469c9d5c0d25d1/states.json"},
  @{Row=47; Date=44299; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8482fd1f83a924b4da62b43afd96f1c7d3828ffe/states.json"},
  @{Row=48; Date=44301; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/ed6993ff4256dbe6b347dcd69c9151e2b60185fb/states.json"},
  @{Row=49; Date=44302; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/05a863ecc79a7568e2cae96170ba08ad957ba885/states.json"},
  @{Row=50; Date=44303; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8b5ae6c27bf498701eb3b33b217714f9f2f9d8c0/states.json"},
  @{Row=51; Date=44304; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/93a4e1f47993b9b81c1a5851dbc6839e3f4707f8/states.json"},
  @{Row=52; Date=44305; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/c73948dd87ed84471e50eeb13e92efe255b943cf/states.json"},
  @{Row=53; Date=44306; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/f194e99b69bbc39ae8387b236c9041aa442f6bb1/states.json"},
  @{Row=54; Date=44307; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/0ba1dd7c6c6eadf91f897c48286751b7b5f2b297/states.json"},
  @{Row=55; Date=44308; Url="https://raw.githubusercontent.com/simonw/cdc-vaccination-history/c7e3d38d73e1819a278db0d363a17957127e0dc8/states.json"}
)

foreach ($e in $entries) {
  $ws.Cells.Item($e.Row, 1).Value = $e.Date
  $ws.Cells.Item($e.Row, 2).Value = $e.Url
}

# Match the saved view state: scrolled so row 37 is at the top, with C56 selected.
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("C56").Select()
